$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task_Table")

$xlLeft = -4131

# ---------------------------------------------------------------------------
# New rows 127-140 for the Cutover execution plan: Dynasys / TMS - Kewill /
# E-Kanban / PricingSpot sections. The shared-string table is append-only,
# so cell values are written in the exact order the original author typed
# them (jumping ahead to type section headers/details, then doubling back
# to fill in a couple of rows) so the regenerated uniqueCount/order matches.
# ---------------------------------------------------------------------------

# Row 127 - Dynasys (section header)
$ws.Cells.Item(127, 2).Value = "Dynasys"
$ws.Cells.Item(127, 2).Font.Bold = $true
$ws.Cells.Item(127, 2).HorizontalAlignment = $xlLeft
$ws.Cells.Item(127, 2).IndentLevel = 1

# Row 128 - Master Data Uploads (sub item)
$ws.Cells.Item(128, 2).Value = "Master Data Uploads"
$ws.Cells.Item(128, 2).IndentLevel = 2

# Row 129 - Server upgrade 2003 to 2008 (sub item)
$ws.Cells.Item(129, 2).Value = "Server upgrade 2003 to 2008"
$ws.Cells.Item(129, 2).IndentLevel = 2

# Row 127, column C - duration
$ws.Cells.Item(127, 3).Value = "2 day"

# Row 130 - TMS - Kewill (section header)
$ws.Cells.Item(130, 2).Value = "TMS - Kewill"
$ws.Cells.Item(130, 2).Font.Bold = $true
$ws.Cells.Item(130, 2).HorizontalAlignment = $xlLeft
$ws.Cells.Item(130, 2).IndentLevel = 1

# Row 133 - E-Kanban (section header) - typed ahead of rows 131/132
$ws.Cells.Item(133, 2).Value = "E-Kanban"
$ws.Cells.Item(133, 2).Font.Bold = $true
$ws.Cells.Item(133, 2).HorizontalAlignment = $xlLeft
$ws.Cells.Item(133, 2).IndentLevel = 1

# Row 134 - Stop application (detail)
$ws.Cells.Item(134, 2).Value = "Stop application"
$ws.Cells.Item(134, 2).HorizontalAlignment = $xlLeft
$ws.Cells.Item(134, 2).IndentLevel = 3

# Row 135 - Change the web.config to SAP FOU (detail)
$ws.Cells.Item(135, 2).Value = "Change the web.config to SAP FOU"
$ws.Cells.Item(135, 2).HorizontalAlignment = $xlLeft
$ws.Cells.Item(135, 2).IndentLevel = 3

# Row 136 - Start application (detail)
$ws.Cells.Item(136, 2).Value = "Start application"
$ws.Cells.Item(136, 2).HorizontalAlignment = $xlLeft
$ws.Cells.Item(136, 2).IndentLevel = 3

# Row 132 - Specify the TMS Kewill WebService URL Endpoint... (sub item, typed before row 131)
$ws.Cells.Item(132, 2).Value = "Specify the TMS Kewill WebService URL Endpoint on SAP Foundation server "
$ws.Cells.Item(132, 2).IndentLevel = 2

# Row 131 - Specify the SAP Foundation URL... (sub item)
$ws.Cells.Item(131, 2).Value = "Specify the SAP Foundation URL on TMS Kewill server to POST inbound shipments "
$ws.Cells.Item(131, 2).IndentLevel = 2

# Row 137 - PricingSpot (section header)
$ws.Cells.Item(137, 2).Value = "PricingSpot "
$ws.Cells.Item(137, 2).Font.Bold = $true
$ws.Cells.Item(137, 2).HorizontalAlignment = $xlLeft
$ws.Cells.Item(137, 2).IndentLevel = 1

# Row 139 - stop windows service (detail) - typed before row 138
$ws.Cells.Item(139, 2).Value = "In the server friawotcssisp1 Stop windows service WSDespesasNutradePRO"
$ws.Cells.Item(139, 2).HorizontalAlignment = $xlLeft
$ws.Cells.Item(139, 2).IndentLevel = 3

# Row 140 - stop IIS site (detail)
$ws.Cells.Item(140, 2).Value = "In the server friawotcssisp1 search for http://syngenta1.pro.intra/pricingspot in the IIS and stop it"
$ws.Cells.Item(140, 2).HorizontalAlignment = $xlLeft
$ws.Cells.Item(140, 2).IndentLevel = 3

# Row 138 - backup database (detail)
$ws.Cells.Item(138, 2).Value = "In the FRIAPSQLGL02\INSGLOB02 instance make the backup for the database PricingSPoT"
$ws.Cells.Item(138, 2).HorizontalAlignment = $xlLeft
$ws.Cells.Item(138, 2).IndentLevel = 3

# Column C durations reusing the existing "1 day" shared string
$ws.Cells.Item(128, 3).Value = "1 day"
$ws.Cells.Item(129, 3).Value = "1 day"

# Column A - literal numeric IDs 171-184
$ids = @{127=171;128=172;129=173;130=174;131=175;132=176;133=177;134=178;135=179;136=180;137=181;138=182;139=183;140=184}
foreach ($r in 127..140) {
    $ws.Cells.Item($r, 1).Value = $ids[$r]
}

# Column B width - widened to fit the new, longer task descriptions.
$ws.Columns.Item(2).ColumnWidth = 92.5

# Sheet view: scroll down to show the newly-added rows, with B141 selected
# as the next empty row.
$ws.Range("B141").Select()
$excel.ActiveWindow.ScrollRow = 127
$excel.ActiveWindow.ScrollColumn = 1
